$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 95
$ws.Range("I2").Value = 243
$ws.Range("J2").Value = 1130
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 293
$ws.Range("M2").Value = 17
$ws.Range("N2").Value = 190
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 20
$ws.Range("S2").Value = 109
$ws.Range("T2").Value = 180
$ws.Range("U2").Value = 18
$ws.Range("V2").Value = 1710
$ws.Range("X2").Value = 1679
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 33
$ws.Range("AA2").Value = 7
